$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G2").Value = "2016-08-23 11:07:33"

$zhcn.Range("H2").Value = "2016-08-23 11:07:28"
$zhcn.Range("K2").Value = "2016-08-23 11:07:44"

$dede.Range("H2").Value = "2016-08-23 11:07:33"
$dede.Range("K2").Value = "2016-08-23 11:07:52"
